# Applies the "In-Game interaction Commands" rewrite:
#  - Paragraph 6 ("Command for interact in item, location and person")
#    becomes "item, location and person" + " In-Game interaction Commands:"
#    (two runs).
#  - Paragraph 7's opening list ("[ Attack, Snack, get, open, use, search,
#    enter") gets the first letter of get/open/use/search/enter capitalised,
#    with each capitalised letter split into its own run, e.g.
#    "[ Attack, Snack, " / "G" / "et, " / "O" / "pen, " / "U" / "se, " /
#    "S" / "earch, " / "E" / "nter".
#
# Word (as exposed by this host) silently re-merges two *adjacent* runs
# once they end up with identical run formatting, so a plain
# Range.Text = "..." rewrite collapses right back into a single run. To
# force a genuine run boundary at a given offset we briefly toggle Bold on
# (which the host always materialises as a distinct run), make the text
# edit, and only afterwards -- once every offset for this paragraph has
# been touched -- flip Bold back off for each of those sub-ranges. Doing
# the "un-bold" pass last (rather than right after each edit) matters:
# unbolding re-normalises neighbouring runs immediately, which -- if done
# before later same-paragraph offsets are processed -- undoes splits that
# haven't been made yet.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Paragraph 6: "Command for interact in item, location and person"
#           -> "item, location and person In-Game interaction Commands:"
# ---------------------------------------------------------------------

$p6 = $d.Paragraphs.Item(6).Range
$p6ContentEnd = $p6.End - 1                      # exclude the pilcrow
$whole6 = $d.Range($p6.Start, $p6ContentEnd)
$whole6.Text = ""

$p6Start = $d.Paragraphs.Item(6).Range.Start
$ins6 = $d.Range($p6Start, $p6Start)
$ins6.InsertAfter("item, location and person In-Game interaction Commands:")

# Split into two runs right before the " In-Game..." suffix.
$splitAt6 = $p6Start + ("item, location and person".Length)
$p6End = $d.Paragraphs.Item(6).Range.End - 1
$tail6 = $d.Range($splitAt6, $p6End)
$tail6.Bold = 1
$tail6b = $d.Range($splitAt6, $p6End)
$tail6b.Bold = 0

# ---------------------------------------------------------------------
# Paragraph 7: "[ Attack, Snack, get, open, use, search, enter..."
#   get -> Get, open -> Open, use -> Use, search -> Search, enter -> Enter
#   each capitalised first letter becomes its own run.
# ---------------------------------------------------------------------

$p7Start = $d.Paragraphs.Item(7).Range.Start

# Offsets (relative to the paragraph start) of the first letter of each
# word, measured against the original text
# "[ Attack, Snack, get, open, use, search, enter".
$offsets = @(17, 22, 28, 33, 41)

foreach ($off in $offsets) {
    $pos = $p7Start + $off
    $r = $d.Range($pos, $pos + 1)
    $letter = $r.Text
    $r.Bold = 1
    $r.Text = $letter.ToUpper()
}

foreach ($off in $offsets) {
    $pos = $p7Start + $off
    $r = $d.Range($pos, $pos + 1)
    $r.Bold = 0
}

Write-Output "Paragraph 6: $($d.Paragraphs.Item(6).Range.Text)"
Write-Output "Paragraph 7: $($d.Paragraphs.Item(7).Range.Text)"
